$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "sending cluster" vocabulary: rename "MuSCs" to "ECs" and
# move "ECs" so that, combined with the new TPM-based numbers below, the
# sending clusters cycle through ECs / Inflammatory-Mac / Neutrophils /
# Resolving-Mac, and the target clusters become Neutrophils / Resolving-Mac.

# New table data (rows 2-9), columns A-T
$data = @(
    @("ECs","Cd86","Ctla4","Neutrophils",1,0.3333333333333333,0.483405,1.450215,0.003094960828876145,0.003094960828876144,1,0.3333333333333333,0.1258873333333333,0.377662,0.9923353626973638,0.9923353626973638,0.06085456637,0.54769109733,0.003071239076656943,0.003071239076656942),
    @("ECs","Cd86","Ctla4","Resolving-Mac",1,0.3333333333333333,0.483405,1.450215,0.003094960828876145,0.003094960828876144,1,0.3333333333333333,0.0009723333333333333,0.002917,0.007664637302636247,0.007664637302636246,0.000470030795,0.004230277155,0.0000237217522192021,0.00002372175221920209),
    @("Inflammatory-Mac","Cd86","Ctla4","Neutrophils",3,1,105.7018663333333,317.105599,0.676747521934545,0.6767475219345449,1,0.3333333333333333,0.1258873333333333,0.377662,0.9923353626973638,0.9923353626973638,13.30652608105978,119.758734729538,0.6715604976334589,0.6715604976334588),
    @("Inflammatory-Mac","Cd86","Ctla4","Resolving-Mac",3,1,105.7018663333333,317.105599,0.676747521934545,0.6767475219345449,1,0.3333333333333333,0.0009723333333333333,0.002917,0.007664637302636247,0.007664637302636246,0.1027774480314444,0.9249970322829999,0.005187024301086155,0.005187024301086154),
    @("Neutrophils","Cd86","Ctla4","Neutrophils",3,1,18.83134933333333,56.494048,0.1205661682058513,0.1205661682058513,1,0.3333333333333333,0.1258873333333333,0.377662,0.9923353626973638,0.9923353626973638,2.370628350641777,21.335655155776,0.1196420722555848,0.1196420722555848),
    @("Neutrophils","Cd86","Ctla4","Resolving-Mac",3,1,18.83134933333333,56.494048,0.1205661682058513,0.1205661682058513,1,0.3333333333333333,0.0009723333333333333,0.002917,0.007664637302636247,0.007664637302636246,0.01831034866844444,0.164793138016,0.0009240959502664843,0.0009240959502664841),
    @("Resolving-Mac","Cd86","Ctla4","Neutrophils",3,1,31.17437066666666,93.523112,0.1995913490307275,0.1995913490307275,1,0.3333333333333333,0.1258873333333333,0.377662,0.9923353626973638,0.9923353626973638,3.924458391571555,35.320125524144,0.1980615537316631,0.1980615537316631),
    @("Resolving-Mac","Cd86","Ctla4","Resolving-Mac",3,1,31.17437066666666,93.523112,0.1995913490307275,0.1995913490307275,1,0.3333333333333333,0.0009723333333333333,0.002917,0.007664637302636247,0.007664637302636246,0.03031187974488889,0.272806917704,0.001529795299064405,0.001529795299064405)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
